# Refresh the coinranking.com snapshot (price / 1h-volume columns, and a
# handful of re-ranked coin rows where Name/Link/Price/Volume all moved).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.465.14'
$ws.Range("E2").Value = '  +0.85%  '

# Row 3
$ws.Range("D3").Value = '1.879.60'
$ws.Range("E3").Value = '  +1.27%  '

# Row 4
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '''0.7182'
$ws.Range("E5").Value = '  +1.68%  '

# Row 6
$ws.Range("D6").Value = '''240.37'
$ws.Range("E6").Value = '  +0.98%  '

# Row 8 (Dogecoin)
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = '''0.07827'
$ws.Range("E8").Value = '  -2.33%  '

# Row 9 (Cardano)
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.3122'
$ws.Range("E9").Value = '  +3.31%  '

# Row 10
$ws.Range("D10").Value = '''25.16'
$ws.Range("E10").Value = '  +7.36%  '

# Row 11
$ws.Range("E11").Value = '  +0.83%  '

# Row 12
$ws.Range("D12").Value = '1.891.81'
$ws.Range("E12").Value = '  +2.52%  '

# Row 13
$ws.Range("D13").Value = '''0.7289'
$ws.Range("E13").Value = '  +3.89%  '

# Row 14
$ws.Range("D14").Value = '''5.292'
$ws.Range("E14").Value = '  +2.12%  '

# Row 15
$ws.Range("D15").Value = '''91.27'
$ws.Range("E15").Value = '  +1.96%  '

# Row 16
$ws.Range("D16").Value = '29.667.53'
$ws.Range("E16").Value = '  +1.99%  '

# Row 17
$ws.Range("D17").Value = '''5.944'
$ws.Range("E17").Value = '  +2.61%  '

# Row 18
$ws.Range("D18").Value = '''248.74'
$ws.Range("E18").Value = '  +4.88%  '

# Row 19
$ws.Range("D19").Value = '''0.000007881'
$ws.Range("E19").Value = '  +0.00%  '

# Row 20
$ws.Range("D20").Value = '''13.31'
$ws.Range("E20").Value = '  +0.76%  '

# Row 21
$ws.Range("D21").Value = '''0.9995'
$ws.Range("E21").Value = '  +0.04%  '

# Row 22
$ws.Range("D22").Value = '''8.004'
$ws.Range("E22").Value = '  +7.45%  '

# Row 23
$ws.Range("D23").Value = '''0.9996'
$ws.Range("E23").Value = '  -0.08%  '

# Row 24
$ws.Range("E24").Value = '  +9.51%  '

# Row 25
$ws.Range("D25").Value = '''163.86'
$ws.Range("E25").Value = '  +0.68%  '

# Row 26
$ws.Range("D26").Value = '''9.056'
$ws.Range("E26").Value = '  +1.80%  '

# Row 27
$ws.Range("D27").Value = '''18.34'
$ws.Range("E27").Value = '  +1.41%  '

# Row 28
$ws.Range("D28").Value = '''1.365'
$ws.Range("E28").Value = '  -3.72%  '

# Row 29
$ws.Range("D29").Value = '''1.485'
$ws.Range("E29").Value = '  +0.48%  '

# Row 30
$ws.Range("D30").Value = '''4.384'
$ws.Range("E30").Value = '  +0.67%  '

# Row 31
$ws.Range("D31").Value = '''4.150'
$ws.Range("E31").Value = '  +3.20%  '

# Row 32
$ws.Range("D32").Value = '''0.05276'
$ws.Range("E32").Value = '  +1.76%  '

# Row 33
$ws.Range("D33").Value = '''1.945'
$ws.Range("E33").Value = '  +1.34%  '

# Row 34
$ws.Range("E34").Value = '  +3.97%  '

# Row 35
$ws.Range("D35").Value = '''0.7238'
$ws.Range("E35").Value = '  +1.64%  '

# Row 36
$ws.Range("D36").Value = '''2.675'
$ws.Range("E36").Value = '  +1.36%  '

# Row 37
$ws.Range("D37").Value = '''0.01864'
$ws.Range("E37").Value = '  +0.81%  '

# Row 38
$ws.Range("D38").Value = '1.241.53'
$ws.Range("E38").Value = '  +9.65%  '

# Row 39
$ws.Range("D39").Value = '''2.725'
$ws.Range("E39").Value = '  +0.10%  '

# Row 40
$ws.Range("D40").Value = '''0.9057'
$ws.Range("E40").Value = '  -2.71%  '

# Row 41
$ws.Range("D41").Value = '''73.77'
$ws.Range("E41").Value = '  +5.27%  '

# Row 42
$ws.Range("D42").Value = '''6.113'
$ws.Range("E42").Value = '  +3.47%  '

# Row 43 (Quant)
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''103.94'
$ws.Range("E43").Value = '  +1.47%  '

# Row 44 (PaxDollar)
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '''0.9999'
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("D45").Value = '''0.5335'
$ws.Range("E45").Value = '  +0.11%  '

# Row 46 (RocketPoolETH)
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.032.85'
$ws.Range("E46").Value = '  +5.25%  '

# Row 47 (RenderToken)
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '''1.767'
$ws.Range("E47").Value = '  +0.56%  '

# Row 48
$ws.Range("D48").Value = '''2.917'
$ws.Range("E48").Value = '  +12.89%  '

# Row 49 (BabyDogeCoin)
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '''0.00000000120'
$ws.Range("E49").Value = '  +0.58%  '

# Row 50
$ws.Range("D50").Value = '''0.4336'
$ws.Range("E50").Value = '  +2.07%  '

# Row 51 (EnergySwap)
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''9.298'
$ws.Range("E51").Value = '  +1.44%  '
